$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.572.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.925.02'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.74'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.62%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4754'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2922'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06785'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '104.87'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +9.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.43'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.914.62'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07729'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.352'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6737'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '287.61'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.612.36'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007639'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.01'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.162.85'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.50%  '
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.485'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +9.88%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.297'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.416'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.20'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.81'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +7.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.135'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +10.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1087'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.363'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.186'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.140'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +9.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05065'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7433'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.160'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02077'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.747'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.066'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '111.33'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8851'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.964'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.13%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4382'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '67.50'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.281'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.388'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.44%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.70'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +15.28%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1234'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.06%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.29'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.31%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4071'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.57%  '
